# Commiting Extent reports update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateSTP")

# Update Data1/Data2 Full Name / Short Name sample values
$ws.Range("F2").Value = "FullName114014"
$ws.Range("F3").Value = "Short114014"
$ws.Range("G2").Value = "FullName214011"
$ws.Range("G3").Value = "Short214011"

# Update Community Organizer / Experts values
$ws.Range("D5").Value = "ukumar1"
$ws.Range("E5").Value = "ukumar1"
$ws.Range("F5").Value = "ukumar1"
$ws.Range("G5").Value = "ukumar1"
$ws.Range("G6").Value = "ukumar1"

# New annotations in column H
$ws.Range("H5").Value = "Uday kumar"
$ws.Range("H6").Value = "Uday Kumar"

# Add a left border to H5 to match the new style introduced
$ws.Range("H5").Borders.Item(7).LineStyle = 1
$ws.Range("H5").Borders.Item(7).Weight = 2

# Move the active cell/selection to C4 as in the authored workbook
$ws.Range("C4").Select()
